# Auto-update draw results: append today's Pick 4 draw as a new row at the
# bottom of the "Results" sheet (row 66), matching the existing table's
# layout (Date, Game, Phase, Result, InsertedAt) and keeping every new cell
# stored as plain text (same as all the pre-existing rows in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 66
$rng = $ws.Range("A$row`:E$row")

# The Date ("2025-11-21") and Phase ("251121") values look like a date /
# a number to Excel's auto-detection, so force text storage on the whole
# new row BEFORE writing any value -- otherwise Excel would silently turn
# them into a date serial / numeric value instead of literal text.
$rng.NumberFormat = "@"

$ws.Range("A$row").Value = "2025-11-21"
$ws.Range("B$row").Value = "Pick 4"
$ws.Range("C$row").Value = "251121"
$ws.Range("D$row").Value = "5-3-0-9"
$ws.Range("E$row").Value = "2025-11-21T21:38:50.505+04:00"

# Put the number format back to the workbook default ("Normal" style) so
# the new row's cells match the rest of the sheet (no explicit text format
# override lingering on them), while the values stay text because they
# were entered while the cells were formatted as text.
$rng.Style = "Normal"
